$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26; existing rows 26-45 shift down to 27-46.
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the latest weekly data point.
# Columns that stay the same as the (old) row 26 / (new) row 27 record:
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 44529
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112009
$ws.Cells.Item(26, 7).Value = "Acelga"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 200
$ws.Cells.Item(26, 11).Value = 1000
$ws.Cells.Item(26, 12).Value = 1200
$ws.Cells.Item(26, 13).Value = 1100
$ws.Cells.Item(26, 14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 367
$ws.Cells.Item(26, 17).Value = 3
$ws.Cells.Item(26, 18).Value = "Hortaliza"
